# Slide 4 ("CVM Instructions") body placeholder "Rectangle 3" gets two of
# its bullet lines reworded for clarification (per commit message).
$p = $ppt.ActivePresentation
$s = $p.Slides.Item(4)
$shape = $s.Shapes.Item("Rectangle 3")
$tf = $shape.TextFrame
$tr = $tf.TextRange

# Paragraph 2: "Some instructions take one or two immediate operands, ..."
# Assigning the final text directly (when it shares a long common substring
# with the existing text) causes the engine to diff the old/new text and
# split the run into several <a:r> elements that each keep the identical
# (matching) or changed (non-matching) segment. Since the real edit is a
# single clean run with the new sentence, first overwrite the paragraph
# with unrelated placeholder text (no overlap with either the old or new
# wording) and then set the real text on top of that -- this keeps the
# paragraph a single run, matching how the slide was actually authored.
$para2 = $tr.Paragraphs(2, 1)
$para2.Text = "PLACEHOLDER_NO_OVERLAP_AAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAAA"
$para2b = $tr.Paragraphs(2, 1)
$para2b.Text = "Some instructions take an immediate operand, which is always located immediately following the instruction in memory."

# Paragraph 3: "Depending on the opcode, an argument can be"
$para3 = $tr.Paragraphs(3, 1)
$para3.Text = "PLACEHOLDER_NO_OVERLAP_BBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBBB"
$para3b = $tr.Paragraphs(3, 1)
$para3b.Text = "Depending on the opcode, an immediate operand can be"
